# PAS_driver_summary.xlsx - Weekly driver report update for 2025-04-21
#
# "Bad Drivers" table (rows 1-5ish): iwlwifi and the AX211-22.250.0.4 rows
# get refreshed counts, the AX201-22.200.2.1 row drops out entirely (no
# longer meets the "bad driver" threshold this week), and the Totals row
# is recomputed and shifts up to fill the gap.
#
# "Good Drivers" table (rows 11+ after the shift): the AX201-23.100.0.4
# row drops out, and the five oldest AX201 rows at the bottom age out of
# the report, leaving only the three AX211 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -------------------------------------------------

# Row 3 "iwlwifi": Client Count 7 -> 6, Good Roaming 91.7 -> 91.8
$ws.Range("B3").Value = 6
$ws.Range("D3").Value = 91.8

# Row 4 "Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4": Client Count 9 -> 8,
# Critical Minutes 639 -> 642
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 642

# Row 5 "Intel(R) Wi-Fi 6 AX201 160MHz - 22.200.2.1" is gone this week -
# delete it outright so the Totals row slides up into row 5.
$ws.Rows(5).Delete()

# Totals row (now row 5): Client Count 23 -> 14, Critical Minutes 5348 -> 4940
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 4940

# --- Good Drivers table --------------------------------------------------

# "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4" (row 14 before the table
# moved up, now row 13 after the Bad Drivers row deletion above) drops off
# the good-driver list this week.
$ws.Rows(13).Delete()

# The five legacy AX201 rows at the bottom of the table (previously rows
# 18-22, now rows 16-20 after the two deletions above) age out, leaving
# only the three AX211 rows.
$ws.Range("A16:A20").EntireRow.Delete()

# The sheet historically carries formatting out through column J / row 20
# (see the <cols> definitions and the original trailing blank rows) even
# though the live data only reaches column E. Touch the bottom-right
# corner cell so the sheet's used range still reaches J20 after the
# trim, matching the refreshed layout.
$ws.Cells.Item(20, 10).Style = "Normal"
